# Apply latest cryptos snapshot values (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.774.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.33%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.797.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.81%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.31%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "433.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.67%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.22%  "

$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("E9").Value = "  -1.12%  "

$ws.Range("E10").Value = "  -11.07%  "

$ws.Range("E11").Value = "  -15.93%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.58"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.417.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.29%  "

$ws.Range("E16").Value = "  -0.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.818.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.75%  "

$ws.Range("E18").Value = "  +1.93%  "

$ws.Range("E19").Value = "  +4.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "66.773.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.97%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "411.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.15%  "

$ws.Range("E23").Value = "  +5.86%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "36.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.69%  "

$ws.Range("E26").Value = "  +5.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +33.70%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.56%  "

$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "719.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.17%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.136"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.60%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +11.01%  "

$ws.Range("E33").Value = "  +0.63%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "41.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.57%  "

$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.150"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.43%  "

$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.66"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +27.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0476"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +38.74%  "

$ws.Range("E41").Value = "  -3.68%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₃0693"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -11.53%  "

$ws.Range("E43").Value = "  +3.70%  "

$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.323"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.05%  "

$ws.Range("E47").Value = "  -0.06%  "

$ws.Range("E48").Value = "  +3.38%  "

$ws.Range("E49").Value = "  -0.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.06%  "
